$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("L28:L29").AutoFill($ws.Range("L28:L32"), 0)
Write-Output $ws.Range("L30").Formula()
Write-Output $ws.Range("L31").Formula()
Write-Output $ws.Range("L32").Formula()
